# The "SearchEngine :" title slide (which only contains the title text and a
# standalone picture) is being dropped from the deck. Its content had already
# been folded into the following "BasicAnalysis :" slide (with the BigQuery
# bullet list), which is why that slide is not touched here: once the
# SearchEngine slide is removed, the BasicAnalysis slide naturally slides up
# to take its place (old slide #13 -> new slide #12), and every following
# slide shifts up by one position as well.

$p = $ppt.ActivePresentation

$slideToRemove = $p.Slides.Item(12)
$slideToRemove.Delete()
